$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Data")
$ws2 = $wb.Worksheets.Item("Metadata")

# --- Data sheet: insert 5 new rows at top for years 2024..2020 ---
$ws1.Range("A2:A6").EntireRow.Insert()

$years = @(2024, 2023, 2022, 2021, 2020)
$vals = @(2.1, 2.5, 2.3, 2.8, 3.2)
for ($i = 0; $i -lt 5; $i++) {
    $r = 2 + $i
    $ws1.Cells.Item($r, 1).Value = [string]$years[$i]
    $ws1.Cells.Item($r, 2).Value = $vals[$i]
}

Write-Output "Data sheet updated"
